$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting existing rows 118-167 down to 119-168.
$ws.Rows("118:118").Insert()

# The newly inserted row 118 is blank; populate it with a copy of what is
# (after the shift) row 119's data -- i.e. the row that used to be row 118 --
# then change only the date (column D) to the new value.
$srcRow = 119
$dstRow = 118
for ($col = 1; $col -le 20; $col++) {
    $srcCell = $ws.Cells.Item($srcRow, $col)
    $dstCell = $ws.Cells.Item($dstRow, $col)
    $dstCell.Value2 = $srcCell.Value2
}
# Column D (4) carries the custom date display format in this sheet.
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat

# New date for the inserted row: 2023-04-11 (serial 45027)
$ws.Cells.Item($dstRow, 4).Value2 = 45027
